# Weekly update: insert a new price record for "Ají" (Americana (o), Primera)
# dated 2023-01-05 (serial 44931) as the new first row of the data block,
# pushing the existing rows (331-358) down to (332-359).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 331, shifting rows down.
$ws.Rows(331).Insert()

# Populate the newly inserted row 331 with the new observation.
$ws.Cells.Item(331, 1).Value = 2
$ws.Cells.Item(331, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(331, 3).Value = "Coquimbo"
$ws.Cells.Item(331, 4).Value = 44931
$ws.Cells.Item(331, 5).Value = 4
$ws.Cells.Item(331, 6).Value = 100112021
$ws.Cells.Item(331, 7).Value = "Ají"
$ws.Cells.Item(331, 8).Value = "Americana (o)"
$ws.Cells.Item(331, 9).Value = "Primera"
$ws.Cells.Item(331, 10).Value = 160
$ws.Cells.Item(331, 11).Value = 12000
$ws.Cells.Item(331, 12).Value = 13000
$ws.Cells.Item(331, 13).Value = 12500
$ws.Cells.Item(331, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(331, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(331, 16).Value = 500
$ws.Cells.Item(331, 17).Value = 25
$ws.Cells.Item(331, 18).Value = "Hortaliza"
